# Scheduled runner update: refresh currentAveragePrice(NQ/HQ)/LevePrice/LeveProfit
# figures across the crafting-job leve tables (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# with newly pulled market-board data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 89.2
$ws.Range("I2").Value = 101.84615
$ws.Range("J2").Value = 83.111115
$ws.Range("K2").Value = 101.84615
$ws.Range("L2").Value = 83.111115
$ws.Range("M2").Value = 11.15385000000001
$ws.Range("N2").Value = -309.111115

$ws.Range("H15").Value = 1505.1
$ws.Range("I15").Value = 1505.1
$ws.Range("K15").Value = 4515.299999999999
$ws.Range("M15").Value = -4346.299999999999

$ws.Range("H33").Value = 2858464.8
$ws.Range("I33").Value = 3334734
$ws.Range("K33").Value = 3334734
$ws.Range("M33").Value = -3334505

$ws.Range("H62").Value = 7814587
$ws.Range("I62").Value = 9616707
$ws.Range("K62").Value = 9616707
$ws.Range("M62").Value = -9616083

$ws.Range("H65").Value = 7814587
$ws.Range("I65").Value = 9616707
$ws.Range("K65").Value = 48083535
$ws.Range("M65").Value = -48080415

$ws.Range("H135").Value = 1640.4193
$ws.Range("I135").Value = 891.4828
$ws.Range("K135").Value = 8023.3452
$ws.Range("M135").Value = -5488.3452

$ws.Range("H137").Value = 4435.3687
$ws.Range("I137").Value = 2374.1904
$ws.Range("K137").Value = 7122.5712
$ws.Range("M137").Value = -4572.5712

$ws.Range("H138").Value = 6556.206
$ws.Range("I138").Value = 4859.125
$ws.Range("K138").Value = 14577.375
$ws.Range("M138").Value = -9437.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 8323.125
$ws.Range("I45").Value = 4146.25
$ws.Range("J45").Value = 12500
$ws.Range("K45").Value = 4146.25
$ws.Range("L45").Value = 12500
$ws.Range("M45").Value = -3769.25
$ws.Range("N45").Value = -13254

$ws.Range("H102").Value = 2559.8
$ws.Range("I102").Value = 3024.75
$ws.Range("K102").Value = 3024.75
$ws.Range("M102").Value = -1402.75

$ws.Range("H132").Value = 4132.5
$ws.Range("I132").Value = 2099.7778
$ws.Range("K132").Value = 6299.3334
$ws.Range("M132").Value = -3769.3334

$ws.Range("H135").Value = 55143.2
$ws.Range("J135").Value = 55143.2
$ws.Range("L135").Value = 55143.2
$ws.Range("N135").Value = -65283.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2935.2
$ws.Range("I20").Value = 2750.5
$ws.Range("K20").Value = 2750.5
$ws.Range("M20").Value = -2503.5

$ws.Range("H99").Value = 3289.087
$ws.Range("I99").Value = 3046.625
$ws.Range("J99").Value = 3843.2856
$ws.Range("K99").Value = 3046.625
$ws.Range("L99").Value = 3843.2856
$ws.Range("M99").Value = -1548.625
$ws.Range("N99").Value = -6839.2856

$ws.Range("H105").Value = 251144.25
$ws.Range("I105").Value = 251144.25
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 251144.25
$ws.Range("L105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("N105").Value = -249397.25

$ws.Range("H107").Value = 1999.4
$ws.Range("I107").Value = 2248.3125
$ws.Range("J107").Value = 1003.75
$ws.Range("K107").Value = 2248.3125
$ws.Range("L107").Value = 1003.75
$ws.Range("M107").Value = -328.3125
$ws.Range("N107").Value = -4843.75

$ws.Range("H114").Value = 60683
$ws.Range("J114").Value = 60683
$ws.Range("L114").Value = 60683
$ws.Range("N114").Value = -69361

$ws.Range("H132").Value = 58333.332
$ws.Range("J132").Value = 58333.332
$ws.Range("L132").Value = 58333.332
$ws.Range("N132").Value = -68453.33199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 298115.03
$ws.Range("I31").Value = 437014.2
$ws.Range("J31").Value = 7689.5454
$ws.Range("K31").Value = 437014.2
$ws.Range("L31").Value = 7689.5454
$ws.Range("M31").Value = -436719.2
$ws.Range("N31").Value = -8279.545399999999

$ws.Range("H34").Value = 298115.03
$ws.Range("I34").Value = 437014.2
$ws.Range("J34").Value = 7689.5454
$ws.Range("K34").Value = 437014.2
$ws.Range("L34").Value = 7689.5454
$ws.Range("M34").Value = -436812.2
$ws.Range("N34").Value = -8093.5454

$ws.Range("H105").Value = 2104.8572
$ws.Range("I105").Value = 1347
$ws.Range("K105").Value = 1347
$ws.Range("M105").Value = 400

$ws.Range("H122").Value = 2568.625
$ws.Range("I122").Value = 1832.8422
$ws.Range("K122").Value = 5498.5266
$ws.Range("M122").Value = -3048.5266

$ws.Range("H132").Value = 4564.7617
$ws.Range("I132").Value = 2724
$ws.Range("K132").Value = 8172
$ws.Range("M132").Value = -5642

$ws.Range("H134").Value = 261186
$ws.Range("I134").Value = 3474.762
$ws.Range("K134").Value = 10424.286
$ws.Range("M134").Value = -7889.286

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 552
$ws.Range("I14").Value = 552
$ws.Range("K14").Value = 1656
$ws.Range("M14").Value = -1483

$ws.Range("H51").Value = 2249
$ws.Range("I51").Value = 3000
$ws.Range("K51").Value = 9000
$ws.Range("M51").Value = -8540

$ws.Range("H55").Value = 6995.5
$ws.Range("I55").Value = 1173.6666
$ws.Range("J55").Value = 11361.875
$ws.Range("K55").Value = 3520.9998
$ws.Range("L55").Value = 34085.625
$ws.Range("M55").Value = -3343.9998
$ws.Range("N55").Value = -34439.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 2312808.2
$ws.Range("J3").Value = 1435785
$ws.Range("L3").Value = 1435785
$ws.Range("N3").Value = -1436017

$ws.Range("H10").Value = 3383334.2
$ws.Range("I10").Value = 5050001.5
$ws.Range("K10").Value = 5050001.5
$ws.Range("M10").Value = -5049832.5

$ws.Range("H43").Value = 26206.5
$ws.Range("I43").Value = 16683
$ws.Range("K43").Value = 16683
$ws.Range("M43").Value = -16532

$ws.Range("H70").Value = 71434490
$ws.Range("I70").Value = 6538.222
$ws.Range("J70").Value = 200004800
$ws.Range("K70").Value = 6538.222
$ws.Range("L70").Value = 200004800
$ws.Range("M70").Value = -6268.222
$ws.Range("N70").Value = -200005340

$ws.Range("H73").Value = 71434490
$ws.Range("I73").Value = 6538.222
$ws.Range("J73").Value = 200004800
$ws.Range("K73").Value = 6538.222
$ws.Range("L73").Value = 200004800
$ws.Range("M73").Value = -5602.222
$ws.Range("N73").Value = -200006672

$ws.Range("H113").Value = 2527502.5
$ws.Range("I113").Value = 3336670
$ws.Range("K113").Value = 3336670
$ws.Range("M113").Value = -3334500

$ws.Range("H122").Value = 2896.0667
$ws.Range("J122").Value = 5979
$ws.Range("L122").Value = 17937
$ws.Range("N122").Value = -22837

$ws.Range("H135").Value = 100068150
$ws.Range("J135").Value = 100068150
$ws.Range("L135").Value = 100068150
$ws.Range("N135").Value = -100078290

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 534019.7
$ws.Range("I7").Value = 11850.5
$ws.Range("J7").Value = 775020.9
$ws.Range("K7").Value = 11850.5
$ws.Range("L7").Value = 775020.9
$ws.Range("M7").Value = -11738.5
$ws.Range("N7").Value = -775244.9

$ws.Range("H46").Value = 4925
$ws.Range("I46").Value = 4925
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 4925
$ws.Range("L46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -4737

$ws.Range("H53").Value = 500
$ws.Range("I53").Value = 500
$ws.Range("K53").Value = 500
$ws.Range("M53").Value = 18

$ws.Range("H61").Value = 4230
$ws.Range("I61").Value = 1766.5714
$ws.Range("J61").Value = 7104
$ws.Range("K61").Value = 1766.5714
$ws.Range("L61").Value = 7104
$ws.Range("M61").Value = -1564.5714
$ws.Range("N61").Value = -7508

$ws.Range("H113").Value = 4230
$ws.Range("I113").Value = 1766.5714
$ws.Range("J113").Value = 7104
$ws.Range("K113").Value = 1766.5714
$ws.Range("L113").Value = 7104
$ws.Range("M113").Value = 403.4286
$ws.Range("N113").Value = -11444

$ws.Range("H122").Value = 912470.0600000001
$ws.Range("I122").Value = 627877.3
$ws.Range("J122").Value = 1671384.1
$ws.Range("K122").Value = 1883631.9
$ws.Range("L122").Value = 5014152.300000001
$ws.Range("M122").Value = -1881181.9
$ws.Range("N122").Value = -5019052.300000001

$ws.Range("H126").Value = 534019.7
$ws.Range("I126").Value = 11850.5
$ws.Range("J126").Value = 775020.9
$ws.Range("K126").Value = 35551.5
$ws.Range("L126").Value = 2325062.7
$ws.Range("M126").Value = -33081.5
$ws.Range("N126").Value = -2330002.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 100000
$ws.Range("J70").Value = 100000
$ws.Range("L70").Value = 100000
$ws.Range("N70").Value = -100630

$ws.Range("H73").Value = 100000
$ws.Range("J73").Value = 100000
$ws.Range("L73").Value = 100000
$ws.Range("N73").Value = -102184

$ws.Range("H113").Value = 1153.5238
$ws.Range("J113").Value = 1592.25
$ws.Range("L113").Value = 4776.75
$ws.Range("N113").Value = -9116.75

$ws.Range("H122").Value = 25002706
$ws.Range("I122").Value = 34484184
$ws.Range("J122").Value = 6077.636
$ws.Range("K122").Value = 103452552
$ws.Range("L122").Value = 18232.908
$ws.Range("M122").Value = -103450102
$ws.Range("N122").Value = -23132.908

$ws.Range("H126").Value = 7166.3335
$ws.Range("I126").Value = 6749.5
$ws.Range("K126").Value = 20248.5
$ws.Range("M126").Value = -17778.5

$ws.Range("H133").Value = 70396
$ws.Range("J133").Value = 70396
$ws.Range("L133").Value = 70396
$ws.Range("N133").Value = -80516

$ws.Range("H135").Value = 48332.5
$ws.Range("J135").Value = 48332.5
$ws.Range("L135").Value = 48332.5
$ws.Range("N135").Value = -58472.5

$ws.Range("H136").Value = 88393.17999999999
$ws.Range("I136").Value = 20976.666
$ws.Range("K136").Value = 62929.99800000001
$ws.Range("M136").Value = -60379.99800000001
